$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_19")

$ws.Range("C5").Value = 3773
$ws.Range("D5").Value = 538603
$ws.Range("E5").Value = 556.41099999999994

$ws.Range("C6").Value = 3749
$ws.Range("D6").Value = 535699
$ws.Range("E6").Value = 552.31799999999998

$ws.Range("C7").Value = 2277
$ws.Range("D7").Value = 326498
$ws.Range("E7").Value = 512.726

$ws.Range("C8").Value = 3147
$ws.Range("D8").Value = 450376
$ws.Range("E8").Value = 534.46900000000005

$ws.Range("C9").Value = 3674
$ws.Range("D9").Value = 523961.99999999994
$ws.Range("E9").Value = 548.15899999999999

$ws.Range("C10").Value = 3784.9999999999995
$ws.Range("D10").Value = 540437
$ws.Range("E10").Value = 565.72699999999998

$ws.Range("B2").Select()
